# Remove the trailing "Ver no Jupiter ..." / "© 2020 ..." footer block
# (and the blank paragraph that separates it from the bibliography entry),
# leaving the bibliography's last line, the existing trailing blank
# paragraph, and the page-break paragraph untouched.

$d = $word.ActiveDocument

# Anchor 1: just after the last character of the final bibliography
# paragraph's text, but before its paragraph mark.
$r1 = $d.Content
$ok1 = $r1.Find.Execute( `
    "Editora Protec, 1991. PROVENZA, F. Projetista de Máquinas . Editora Protec, 1991.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok1) {
    throw "Could not locate the bibliography anchor text."
}
# Step past that paragraph's end-of-paragraph mark so the delete range
# starts at the beginning of the following (blank) paragraph.
$r1.MoveEnd(1, 1) | Out-Null
$startDelete = $r1.End

# Anchor 2: end of the copyright/footer paragraph, including its
# paragraph mark, so the whole paragraph is removed.
$r2 = $d.Content
$ok2 = $r2.Find.Execute( `
    "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok2) {
    throw "Could not locate the footer/copyright anchor text."
}
$r2.MoveEnd(1, 1) | Out-Null
$endDelete = $r2.End

# Delete the blank paragraph + "Ver no Jupiter ..." paragraph +
# "© 2020 ..." paragraph in one shot.
$deleteRange = $d.Range($startDelete, $endDelete)
$deleteRange.Delete()
